$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.370830059051514
$ws.Range("B1").Value = 1.95443594455719
$ws.Range("C1").Value = 2.95820689201355
$ws.Range("D1").Value = 3.71523380279541
$ws.Range("E1").Value = 1.005772590637207
